# Apply edits to add a new "PO Forecast" sheet with PO forecast data,
# and rename the "Requested quantity" header on the two existing sheets.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Rename headers on existing sheets.
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet so it lands
# at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# Header row.
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

$ws.Range("A2").Value = 44976.99999999999
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = -227.5984874810388
$ws.Range("D2").Value = 130.1972359667564
$ws.Range("A3").Value = 44997.99999999999
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = -175.8766190904722
$ws.Range("D3").Value = 177.6068544014902
$ws.Range("A4").Value = 45004.99999999999
$ws.Range("B4").Value = 16
$ws.Range("C4").Value = -174.6551256129793
$ws.Range("D4").Value = 190.3566237306973
$ws.Range("A5").Value = 45025.99999999999
$ws.Range("B5").Value = 70
$ws.Range("C5").Value = -121.7808823939456
$ws.Range("D5").Value = 237.9593099729635
$ws.Range("A6").Value = 45039.99999999999
$ws.Range("B6").Value = 106
$ws.Range("C6").Value = -70.62482955813765
$ws.Range("D6").Value = 293.5413726314707
$ws.Range("A7").Value = 45046.99999999999
$ws.Range("B7").Value = 124
$ws.Range("C7").Value = -50.69268524986472
$ws.Range("D7").Value = 309.4984397701356
$ws.Range("A8").Value = 45053.99999999999
$ws.Range("B8").Value = 141
$ws.Range("C8").Value = -27.8913574960773
$ws.Range("D8").Value = 311.2793341548079
$ws.Range("A9").Value = 45060.99999999999
$ws.Range("B9").Value = 159
$ws.Range("C9").Value = -16.2723345994354
$ws.Range("D9").Value = 345.661499011789
$ws.Range("A10").Value = 45067.99999999999
$ws.Range("B10").Value = 177
$ws.Range("C10").Value = -1.192054104215695
$ws.Range("D10").Value = 363.0556993574307
$ws.Range("A11").Value = 45074.99999999999
$ws.Range("B11").Value = 195
$ws.Range("C11").Value = 3.843102694654758
$ws.Range("D11").Value = 352.6256052415297
$ws.Range("A12").Value = 45088.99999999999
$ws.Range("B12").Value = 231
$ws.Range("C12").Value = 53.18746726605176
$ws.Range("D12").Value = 411.3431517059602
$ws.Range("A13").Value = 45137.99999999999
$ws.Range("B13").Value = 356
$ws.Range("C13").Value = 173.4324792247165
$ws.Range("D13").Value = 532.6383248537745
$ws.Range("A14").Value = 45144.99999999999
$ws.Range("B14").Value = 374
$ws.Range("C14").Value = 202.5447311632138
$ws.Range("D14").Value = 544.910370206228
$ws.Range("A15").Value = 45151.99999999999
$ws.Range("B15").Value = 391
$ws.Range("C15").Value = 212.0357201503346
$ws.Range("D15").Value = 571.8255208374843
$ws.Range("A16").Value = 45158.99999999999
$ws.Range("B16").Value = 409
$ws.Range("C16").Value = 243.3097947670885
$ws.Range("D16").Value = 599.3442023275821
$ws.Range("A17").Value = 45165.99999999999
$ws.Range("B17").Value = 427
$ws.Range("C17").Value = 246.7917240841895
$ws.Range("D17").Value = 601.1755001419048
$ws.Range("A18").Value = 45172.99999999999
$ws.Range("B18").Value = 445
$ws.Range("C18").Value = 258.0437501490431
$ws.Range("D18").Value = 613.2028270265971
$ws.Range("A19").Value = 45179.99999999999
$ws.Range("B19").Value = 463
$ws.Range("C19").Value = 278.3308085840644
$ws.Range("D19").Value = 627.6427589274093
$ws.Range("A20").Value = 45186.99999999999
$ws.Range("B20").Value = 481
$ws.Range("C20").Value = 307.7564605343632
$ws.Range("D20").Value = 653.0082809540017
$ws.Range("A21").Value = 45193.99999999999
$ws.Range("B21").Value = 499
$ws.Range("C21").Value = 314.6682733426488
$ws.Range("D21").Value = 682.3178836271785

# Match formatting of the other sheets: bold/centered/bordered header row,
# and the date number format on column A.
$wsWeekly.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("C1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$ws.Range("A2:A21").PasteSpecial(-4122)
